# Apply odds-update edits to "Sheet1" of the FlashScore weekly games workbook.
# Each statement below writes one corrected numeric value into the cell
# identified by its A1 reference, matching the author's commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 - Twente vs PSV
$ws.Range("N3").Value = 1.5
$ws.Range("O3").Value = 2.5

# Row 10 - Club Brugge KV vs Royale Union SG
$ws.Range("K10").Value = 13
$ws.Range("AA10").Value = 6.5
$ws.Range("AD10").Value = 201
$ws.Range("AI10").Value = 34

# Row 11 - Nueve de Octubre vs Cumbaya
$ws.Range("G11").Value = 2.8
$ws.Range("I11").Value = 2.47
$ws.Range("R11").Value = 2
$ws.Range("T11").Value = 6.9
$ws.Range("U11").Value = 12.5
$ws.Range("W11").Value = 32
$ws.Range("X11").Value = 29
$ws.Range("AB11").Value = 18

# Row 12 - Chacaritas vs Guayaquil City
$ws.Range("H12").Value = 3.2
$ws.Range("I12").Value = 1.95
$ws.Range("W12").Value = 60
$ws.Range("AA12").Value = 6.4
$ws.Range("AF12").Value = 7.9
$ws.Range("AH12").Value = 16.5

# Row 16 - Torpedo Kutaisi vs Gagra
$ws.Range("M16").Value = 3.3
$ws.Range("R16").Value = 1.9
$ws.Range("S16").Value = 1.81

# Row 17 - Vestmannaeyjar vs Fram
$ws.Range("H17").Value = 3.7
$ws.Range("L17").Value = 1.2
$ws.Range("M17").Value = 4.1
$ws.Range("O17").Value = 2.2
$ws.Range("S17").Value = 2.32
$ws.Range("V17").Value = 9.5
$ws.Range("AA17").Value = 7.5
$ws.Range("AE17").Value = 11
$ws.Range("AI17").Value = 18

# Row 20 - Auda vs Grobina
$ws.Range("R20").Value = 2.01
$ws.Range("S20").Value = 1.71

# Row 29 - Balestier Khalsa vs Hougang
$ws.Range("G29").Value = 1.85
$ws.Range("I29").Value = 2.7
$ws.Range("M29").Value = 17
$ws.Range("N29").Value = 1.11
$ws.Range("O29").Value = 6.5
$ws.Range("P29").Value = 1.1
$ws.Range("Q29").Value = 7
$ws.Range("R29").Value = 1.17
$ws.Range("S29").Value = 5
$ws.Range("T29").Value = 34
$ws.Range("U29").Value = 26
$ws.Range("V29").Value = 13
$ws.Range("W29").Value = 29
$ws.Range("X29").Value = 15
$ws.Range("AA29").Value = 19
$ws.Range("AB29").Value = 13
$ws.Range("AE29").Value = 41
$ws.Range("AJ29").Value = 15

# Row 30 - Djurgarden vs Oster
$ws.Range("N30").Value = 1.83
$ws.Range("O30").Value = 2.03

# Row 33 - Varnamo vs AIK
$ws.Range("G33").Value = 3.1
$ws.Range("I33").Value = 2.3
$ws.Range("J33").Value = 1.1
$ws.Range("K33").Value = 7
$ws.Range("T33").Value = 7
$ws.Range("W33").Value = 34
$ws.Range("X33").Value = 29
$ws.Range("AF33").Value = 9.5
$ws.Range("AI33").Value = 23
